# Apply updated First_Noticeable_Increase_Index (C), First_Noticeable_Increase_Cumulative_Value (E)
# and Pulse_Width (G) values to each of the Step3_DataPts_* sheets. These values reflect the
# addition of a configurable zero_before_threshold parameter that can zero out dims before the
# noise threshold / First Rise Point.

$wb = $excel.ActiveWorkbook

# Values that are identical across all four Step3_DataPts_* sheets (rows 2-6), keyed by row.
$commonByRow = @{
    2 = @{ C = 87; E = 0.02601518943497679 }
    3 = @{ C = 90; E = 0.01101021647796185 }
    4 = @{ C = 87; E = 0.05151036632781263 }
    5 = @{ C = 87; E = 0.01622240674782391 }
    6 = @{ C = 87; E = 0.03079592233537105 }
}

# Pulse_Width (column G) values differ per sheet/threshold, keyed by sheet name then row.
$gByThresholdSheet = @{
    "Step3_DataPts_0.5" = @{ 2 = 45; 3 = 49; 4 = 35; 5 = 29; 6 = 21 }
    "Step3_DataPts_0.7" = @{ 2 = 62; 3 = 66; 4 = 58; 5 = 58; 6 = 56 }
    "Step3_DataPts_0.8" = @{ 2 = 76; 3 = 76; 4 = 72; 5 = 68; 6 = 67 }
    "Step3_DataPts_0.9" = @{ 2 = 96; 3 = 101; 4 = 96; 5 = 94; 6 = 85 }
}

foreach ($sheetName in $gByThresholdSheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $gByRow = $gByThresholdSheet[$sheetName]

    foreach ($row in 2..6) {
        $ws.Range("C$row").Value = $commonByRow[$row].C
        $ws.Range("E$row").Value = $commonByRow[$row].E
        $ws.Range("G$row").Value = $gByRow[$row]
    }
}
